$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defect List")
$scratch = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1) Create the two new cell styles (cellXfs #30 and #31) that the target
#    workbook needs, by copying an existing "wrap + colored" style and then
#    recoloring the fill to yellow (same fill used by severity "3"), using a
#    scratch cell on the empty Sheet2 so the visible sheet is untouched.
#    - style 30: like D59's old style (fontId3/fillId4/wrap) but fillId2
#    - style 31: like D51's old style (fontId5/fillId4/wrap) but fillId2
# ---------------------------------------------------------------------------
$ws.Range("D59").Copy()
$scratch.Range("A1").PasteSpecial(-4122)          # xlPasteFormats
$scratch.Range("A1").Interior.Color = 65535       # RGB yellow -> existing fillId2
$scratch.Range("A1").Clear()

$ws.Range("D51").Copy()
$scratch.Range("A1").PasteSpecial(-4122)          # xlPasteFormats
$scratch.Range("A1").Interior.Color = 65535       # RGB yellow -> existing fillId2
$scratch.Range("A1").Clear()

# ---------------------------------------------------------------------------
# 2) Severity edits: rows 48, 51 and 59 go from severity 4 (blue) to
#    severity 3 (yellow), reusing the matching yellow-filled style.
# ---------------------------------------------------------------------------
$ws.Range("D48").Value = 3
$ws.Range("D48").Interior.Color = 65535

$ws.Range("D51").Value = 3
$ws.Range("D51").Interior.Color = 65535

$ws.Range("D59").Value = 3
$ws.Range("D59").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3) State edits: row 57 moves from "Active " to "Fixed"; row 60 moves from
#    "Fixed" to the new "Resolved" state.
# ---------------------------------------------------------------------------
$ws.Range("I57").Value = "Fixed"
$ws.Range("I60").Value = "Resolved"

# ---------------------------------------------------------------------------
# 4) New defect row 54 (sheet row 61), copying each cell's format from an
#    existing cell that already carries the right style, then filling in the
#    row's values.
# ---------------------------------------------------------------------------
$ws.Range("A49").Copy()
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("B57").Copy()
$ws.Range("B61").PasteSpecial(-4122)
$ws.Range("C48").Copy()
$ws.Range("C61").PasteSpecial(-4122)
$ws.Range("D57").Copy()
$ws.Range("D61").PasteSpecial(-4122)
$ws.Range("E57").Copy()
$ws.Range("E61").PasteSpecial(-4122)
$ws.Range("F57").Copy()
$ws.Range("F61").PasteSpecial(-4122)
$ws.Range("G57").Copy()
$ws.Range("G61").PasteSpecial(-4122)
$ws.Range("H51").Copy()
$ws.Range("H61").PasteSpecial(-4122)
$ws.Range("I57").Copy()
$ws.Range("I61").PasteSpecial(-4122)

$ws.Range("A61").Value = 54
$ws.Range("B61").Value = "Jennifer"
$ws.Range("C61").Value = 42104
$ws.Range("D61").Value = 2
$ws.Range("E61").Value = "ST 5.0, ST. 7.0 & ST. 8.0"
$ws.Range("F61").Value = "Creating a task"
$ws.Range("G61").Value = "When a task is created with more than two words, the action of submitting a task will not update"
$ws.Range("H61").Value = "Create a task with two or more words "
$ws.Range("I61").Value = "Active "

$ws.Rows.Item(61).RowHeight = 63

# ---------------------------------------------------------------------------
# 5) Scroll position / selection bookkeeping to match the refreshed view.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 56
$ws.Range("E61").Select()

Write-Output "edit complete"
